# Swap the "B" and "C" quarter rows within each year-group (rows 2-5, 6-9, ...,
# i.e. row pairs (3,4) (7,8) (11,12) ... (63,64)) and drop the now-redundant
# F/G columns (家用吸尘器产销率 / 家用吸尘器销售量), which duplicated B/E with a
# shifted row alignment.
#
# Notes on this runtime's COM surface (discovered by experimentation):
#  - `.Value` does not resolve correctly when read back out of a Range/Cell
#    (it yields the COM member-descriptor string instead of the cell's
#    contents); `.Value2` reads correctly, so it is used for every read here.
#  - Writing an empty string via `.Value2 = ""` (or `.ClearContents()`)
#    always normalises the cell to a true blank rather than preserving an
#    empty-string text cell. To avoid needlessly losing the original
#    empty-string typing on cells that don't actually need to change, a
#    cell pair is only written when its two values actually differ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(3,4), @(7,8), @(11,12), @(15,16), @(19,20), @(23,24), @(27,28), @(31,32),
    @(35,36), @(39,40), @(43,44), @(47,48), @(51,52), @(55,56), @(59,60), @(63,64)
)

$scratchRow = 1000

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($c = 1; $c -le 5; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        if ($v1 -ne $v2) {
            $scratch = $ws.Cells.Item($scratchRow, $c)
            $scratch.Value2 = $v1
            $cell1.Value2 = $v2
            $cell2.Value2 = $scratch.Value2
            $scratch.Value2 = ""
        }
    }
}

# Remove the duplicate "产销率" / "销售量" columns (F and G)
$ws.Range("F1:G65").EntireColumn.Delete()
